$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Autonomous" ---
$ws1 = $wb.Worksheets.Item(1)

$senders1 = @{
  1  = "Sender(s): Dynamics"
  12 = "Sender(s): ACU"
  23 = "Sender(s): RES"
  27 = "Sender(s): ACU"
  32 = "Sender(s): ACU"
  41 = "Sender(s): Jetson"
}

foreach ($r in $senders1.Keys) {
  # Copy formatting (bold text + header fill + border) from the adjacent "ID:" cell
  # in column B onto the new column C cell, then set its text.
  $ws1.Range("B$r").Copy()
  $ws1.Range("C$r").PasteSpecial(-4122)
  $ws1.Range("C$r").Value2 = $senders1[$r]
}

# Widen column C from 15 to 21 characters (Excel's ColumnWidth property is
# offset by 5/6 of a character from the stored OOXML column width).
$ws1.Columns.Item(3).ColumnWidth = 21 - (5/6)

# --- Sheet 2: "Autonomous_temporary" ---
$ws2 = $wb.Worksheets.Item(2)

$senders2 = @{
  1  = "Sender(s): ACU"
  5  = "Sender(s): JETSON"
  9  = "Sender(s): VCU"
  13 = "Sender(s): JETSON"
  17 = "Sender(s): ACU"
  22 = "Sender(s): ACU"
  26 = "Sender(s): JETSON"
  30 = "Sender(s): VCU"
  34 = "Sender(s): RES"
}

foreach ($r in $senders2.Keys) {
  $ws2.Range("B$r").Copy()
  $ws2.Range("C$r").PasteSpecial(-4122)
  $ws2.Range("C$r").Value2 = $senders2[$r]
}

# Widen column C from 15 to 19 characters
$ws2.Columns.Item(3).ColumnWidth = 19 - (5/6)

$excel.CutCopyMode = 0
